$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '30.287.93'
    'D3' = '1.934.19'
    'E3' = '  -3.26%  '
    'D4' = '1.001'
    'E4' = '  +0.25%  '
    'D5' = '247.55'
    'E5' = '  -2.99%  '
    'D6' = '0.7178'
    'E6' = '  -11.56%  '
    'D7' = '1.000'
    'D8' = '0.3282'
    'E8' = '  -7.80%  '
    'D9' = '26.65'
    'E9' = '  +2.92%  '
    'D10' = '0.06837'
    'E10' = '  -2.70%  '
    'D11' = '0.8082'
    'D12' = '0.07982'
    'E12' = '  -1.95%  '
    'D13' = '1.932.77'
    'E13' = '  -3.31%  '
    'D14' = '5.443'
    'E14' = '  -1.55%  '
    'D15' = '94.90'
    'E15' = '  -6.51%  '
    'D16' = '14.61'
    'E16' = '  +4.04%  '
    'D17' = '263.96'
    'E17' = '  -3.32%  '
    'D18' = '30.293.24'
    'E18' = '  -3.48%  '
    'B19' = 'Uniswap'
    'C19' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D19' = '5.848'
    'E19' = '  +0.36%  '
    'B20' = 'ShibaInu'
    'C20' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D20' = '0.000007965'
    'E20' = '  +0.01%  '
    'D21' = '2.187.43'
    'E21' = '  -2.86%  '
    'D22' = '0.9997'
    'E22' = '  +0.26%  '
    'D23' = '1.001'
    'E23' = '  +0.27%  '
    'D24' = '6.931'
    'E24' = '  -1.50%  '
    'D25' = '9.769'
    'E25' = '  -1.08%  '
    'D26' = '160.33'
    'E26' = '  -2.50%  '
    'D27' = '2.335'
    'E27' = '  +2.70%  '
    'D28' = '0.1343'
    'E28' = '  -11.85%  '
    'D29' = '19.06'
    'E29' = '  -5.40%  '
    'E30' = '  +0.50%  '
    'D31' = '1.555'
    'E31' = '  -1.30%  '
    'D32' = '4.424'
    'E32' = '  -4.13%  '
    'D33' = '4.240'
    'E33' = '  -2.89%  '
    'D34' = '0.05099'
    'E34' = '  -2.36%  '
    'D35' = '1.209'
    'E35' = '  -0.89%  '
    'D36' = '0.7486'
    'E36' = '  -2.07%  '
    'D37' = '2.740'
    'E37' = '  -0.46%  '
    'D38' = '0.01948'
    'E38' = '  -3.31%  '
    'D39' = '2.822'
    'E39' = '  -3.39%  '
    'D40' = '81.39'
    'E40' = '  +3.37%  '
    'D41' = '6.591'
    'E41' = '  -0.90%  '
    'D42' = '0.4492'
    'E42' = '  -5.69%  '
    'D43' = '2.022'
    'E43' = '  -5.34%  '
    'E44' = '  +0.28%  '
    'D45' = '0.8375'
    'E45' = '  -2.86%  '
    'D46' = '102.49'
    'E46' = '  -2.09%  '
    'D47' = '9.849'
    'E47' = '  -1.96%  '
    'D48' = '7.351'
    'E48' = '  -2.57%  '
    'D49' = '36.35'
    'E49' = '  -1.51%  '
    'D50' = '1.501'
    'E50' = '  +2.80%  '
    'D51' = '0.4128'
}

foreach ($addr in $updates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$addr]
}
